# Add a new bull "AMERICAN FLORIAN" (KI-code 361226) to the
# "Fokstieren stierenkaart" list, right after "SURPRISE rf" (row 16),
# leaving a blank spacer row beneath it before the rest of the list
# continues - matching the author's manual row-insert edit.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Fokstieren stierenkaart
$ws2 = $wb.Worksheets.Item(2)   # Genomics stierenkaart

# Insert two blank rows at row 17 (shifts the old row 17.. down to 19..)
[void]$ws1.Rows("17:18").Insert()

# Fill the first of the two new rows with the new bull's data
$ws1.Range("A17").Value = 361226
$ws1.Range("B17").Value = "AMERICAN FLORIAN"
# row 18 stays blank, as in the source edit

# Restore the on-screen selections to match the saved state
[void]$ws2.Rows("21:21").Select()
[void]$ws1.Activate()
[void]$ws1.Range("E27").Select()
